# edit.ps1 -- apply the "design changes report changes" revision to the
# Technical Specification (Vapour Absorption Chiller) document.
#
# Summary of changes (see the commit diff this implements):
#   1. Remove the leading centered paragraph that only contains the
#      decorative VML picture (w:pict / v:shape, rId7).
#   2. Update a handful of table-cell values (version/date/spec numbers).
#   3. Rename "Evaporator tube material" -> "Evaporator" and its value
#      "Copper" -> "name".
#   4. Flip "NonStandard" -> "Standard" for the Low Temperature Heat
#      exchanger Type row.
#   5. Rework the "Caption Notes" list: drop notes 1 & 2, keep the rest,
#      renumbering them 1-4 with updated wording.

$d = $word.ActiveDocument

function Replace-FirstText($find, $replace, $wholeWord) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($find, $true, $wholeWord, $false, $false, $false, `
                             $true, 1, $false, $replace, 1)
    if (-not $ok) {
        Write-Output "WARNING: text not found -> $find"
    }
}

function Remove-ParagraphByText($target) {
    $n = $d.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $target) {
            $p.Range.Delete()
            return
        }
    }
    Write-Output "WARNING: paragraph not found -> $target"
}

# --- 1. Drop the lone picture paragraph at the very top of the document ---
$firstPara = $d.Paragraphs.Item(1)
if ($firstPara.Range.Text.TrimEnd([char]13, [char]7) -eq "") {
    $firstPara.Range.Delete()
}

# --- 2. Header/info table value updates -----------------------------------
Replace-FirstText "5.1.2.0" "1.0" $true
Replace-FirstText "04/05/2021, 11:23 AM" "05/21/2021, 04:14 PM" $true

# --- 3. Spec table value updates -------------------------------------------
Replace-FirstText "4+4" "2+2" $true
Replace-FirstText "6.8" "2.3" $true
Replace-FirstText "0.00002" "standard" $true
Replace-FirstText "32.6" "33" $true
Replace-FirstText "35" "34.9" $true
Replace-FirstText "4.6" "4.7" $true
Replace-FirstText "0.00005" "standard" $true
Replace-FirstText "3.0" "3.2" $true

# --- 4. Evaporator tube material row ---------------------------------------
Replace-FirstText "Evaporator tube material" "Evaporator" $false
Replace-FirstText "Copper" "name" $true

# --- 5. Low Temperature Heat exchanger Type row -----------------------------
Replace-FirstText "NonStandard" "Standard" $true

# --- 6. Caption Notes list: remove notes 1 & 2, renumber the remainder -----
Remove-ParagraphByText "1. This is an ARI selection"
Remove-ParagraphByText "2. Note----Higher Size Low temp Heat Exchanger Required, Pls refer to Engg. for extra cost"

Replace-FirstText "3. This selection is valid for insulated chiller only." "1. This selection is valid for insulated chiller only." $false
Replace-FirstText "4. For non-insulated chiller, the Capacity and Heat source consumption will vary." "2. For non-insulated chiller, the Capacity and Heat source consumption will vary." $false
Replace-FirstText "5. Plant Room Temperature should be from +5 deg C to +45 deg C" "3. Plant Room Temperature should be from +5 deg C to +45 deg C" $false
Replace-FirstText "6. Please contact Thermax representative / Office for customised specifications." "4. Please contact Thermax representative / Office for customised specifications." $false

Write-Output "done"
